# Update Tnf-Tnfrsf1b LR-pair sheet with new TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 5.599488
$ws.Range("H2").Value = 16.798464
$ws.Range("I2").Value = 0.8205563069596913
$ws.Range("J2").Value = 0.8205563069596913
$ws.Range("M2").Value = 5.916202333333334
$ws.Range("N2").Value = 17.748607
$ws.Range("O2").Value = 0.3515586392055965
$ws.Range("P2").Value = 0.3515586392055965
$ws.Range("Q2").Value = 33.127703971072
$ws.Range("R2").Value = 298.149335739648
$ws.Range("S2").Value = 0.2884736586663189
$ws.Range("T2").Value = 0.2884736586663189

# Row 3
$ws.Range("G3").Value = 5.599488
$ws.Range("H3").Value = 16.798464
$ws.Range("I3").Value = 0.8205563069596913
$ws.Range("J3").Value = 0.8205563069596913
$ws.Range("O3").Value = 0.6159539016771971
$ws.Range("P3").Value = 0.6159539016771971
$ws.Range("Q3").Value = 58.041920291584
$ws.Range("R3").Value = 522.3772826242559
$ws.Range("S3").Value = 0.5054248588176536
$ws.Range("T3").Value = 0.5054248588176536

# Row 4
$ws.Range("G4").Value = 5.599488
$ws.Range("H4").Value = 16.798464
$ws.Range("I4").Value = 0.8205563069596913
$ws.Range("J4").Value = 0.8205563069596913
$ws.Range("M4").Value = 0.5467150000000001
$ws.Range("O4").Value = 0.03248745911720639
$ws.Range("P4").Value = 0.03248745911720639
$ws.Range("Q4").Value = 3.06132408192
$ws.Range("R4").Value = 27.55191673728
$ws.Range("S4").Value = 0.02665778947571883
$ws.Range("T4").Value = 0.02665778947571883

# Row 5
$ws.Range("I5").Value = 0.1794436930403087
$ws.Range("J5").Value = 0.1794436930403087
$ws.Range("M5").Value = 5.916202333333334
$ws.Range("N5").Value = 17.748607
$ws.Range("O5").Value = 0.3515586392055965
$ws.Range("P5").Value = 0.3515586392055965
$ws.Range("Q5").Value = 7.244545550494778
$ws.Range("R5").Value = 65.200909954453
$ws.Range("S5").Value = 0.0630849805392777
$ws.Range("T5").Value = 0.0630849805392777

# Row 6
$ws.Range("I6").Value = 0.1794436930403087
$ws.Range("J6").Value = 0.1794436930403087
$ws.Range("O6").Value = 0.6159539016771971
$ws.Range("P6").Value = 0.6159539016771971
$ws.Range("S6").Value = 0.1105290428595434
$ws.Range("T6").Value = 0.1105290428595434

# Row 7
$ws.Range("I7").Value = 0.1794436930403087
$ws.Range("J7").Value = 0.1794436930403087
$ws.Range("M7").Value = 0.5467150000000001
$ws.Range("O7").Value = 0.03248745911720639
$ws.Range("P7").Value = 0.03248745911720639
$ws.Range("Q7").Value = 0.6694669143283334
$ws.Range("R7").Value = 6.025202228955001
$ws.Range("S7").Value = 0.005829669641487562
$ws.Range("T7").Value = 0.005829669641487562
